$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.036.69"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.14"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9972"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.68"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6325"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9993"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07513"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2941"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.05"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.833.23"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.999"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6714"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.11"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009567"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +5.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.081"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.057.56"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.59"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "226.68"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9984"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.172"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9992"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.16"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1409"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +3.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.540"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +1.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.93"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.501"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.139"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.066"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.197"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05394"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +3.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.857"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7440"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +1.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.140"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.654"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.244.49"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -3.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.764"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.663"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +4.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9039"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9997"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.92"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.977.86"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000123"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +3.22%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.96"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +2.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5099"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4065"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.964"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.762"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +1.24%  "
